$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.52'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '25.36'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.107'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05577'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.501'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.020'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8174'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8466'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06953'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03166'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02891'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09387'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001525'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006194'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.009794'
$ws.Range("E19").Value = '18OneONEBestin24h'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3179'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.748'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04735'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001246'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004638'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009704'
$ws.Range("E27").Value = '26NitroExNTX'
$ws.Range("E28").Value = '27UpBotsUBXT'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03670'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006172'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1052'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002586'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008307'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005302'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1501'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002123'
